$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 7) to the users_engagement sheet, mirroring
# the structure/columns of the existing rows:
#   user_id | username | level | last_message_date | last_response |
#   response_status | level_3_ai_response | subscription_checked |
#   level_4_reminder_sent | decision | notes
$row = 7

$ws.Cells.Item($row, 1).Value  = 6698418542
$ws.Cells.Item($row, 2).Value  = "sourabratabose"
$ws.Cells.Item($row, 3).Value  = -1
$ws.Cells.Item($row, 6).Value  = "unreachable"
$ws.Cells.Item($row, 8).Value  = $false
$ws.Cells.Item($row, 9).Value  = $false
$ws.Cells.Item($row, 11).Value = "Added during extraction"
